$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "22.454.11"
$r.ClearFormats()
$ws.Range("E2").Value = "  +0.14%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.573.11"
$r.ClearFormats()
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.03%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "291.08"
$r.ClearFormats()
$ws.Range("E6").Value = "  +0.00%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.3740"
$r.ClearFormats()
$ws.Range("E7").Value = "  -0.72%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "49.99"
$r.ClearFormats()
$ws.Range("E8").Value = "  +0.02%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.3404"
$r.ClearFormats()
$ws.Range("E9").Value = "  -0.50%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.07558"
$r.ClearFormats()
$ws.Range("E10").Value = "  -1.47%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "1.141"
$r.ClearFormats()
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("E12").Value = "  -0.04%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "21.35"
$r.ClearFormats()
$ws.Range("E13").Value = "  +0.38%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "5.986"
$r.ClearFormats()
$ws.Range("E14").Value = "  -0.08%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "6.952"
$r.ClearFormats()
$ws.Range("E15").Value = "  +0.38%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "1.573.98"
$r.ClearFormats()
$ws.Range("E16").Value = "  -0.13%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.00001124"
$r.ClearFormats()
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("E18").Value = "  +0.62%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.06731"
$r.ClearFormats()
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  -0.01%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "6.269"
$r.ClearFormats()
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("E23").Value = "  +1.18%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "22.449.84"
$r.ClearFormats()
$ws.Range("E24").Value = "  +0.14%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.346"
$r.ClearFormats()
$ws.Range("E25").Value = "  -3.18%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "2.577"
$r.ClearFormats()
$ws.Range("E26").Value = "  -6.08%  "
$ws.Range("E27").Value = "  -0.95%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "148.68"
$r.ClearFormats()
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").Value = "  -0.84%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "125.87"
$r.ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.749.84"
$r.ClearFormats()
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +3.58%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "6.127"
$r.ClearFormats()
$ws.Range("E33").Value = "  -1.55%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "1.977"
$r.ClearFormats()
$ws.Range("E34").Value = "  -1.84%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "9.850"
$r.ClearFormats()
$ws.Range("E35").Value = "  -1.74%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.08424"
$r.ClearFormats()
$ws.Range("E36").Value = "  -1.74%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.377"
$r.ClearFormats()
$ws.Range("E37").Value = "  +1.81%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.02459"
$r.ClearFormats()
$ws.Range("E38").Value = "  -3.96%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.2295"
$r.ClearFormats()
$ws.Range("E39").Value = "  -0.95%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.06529"
$r.ClearFormats()
$ws.Range("E40").Value = "  -0.30%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "5.480"
$r.ClearFormats()
$ws.Range("E41").Value = "  +0.19%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "11.35"
$r.ClearFormats()
$ws.Range("E42").Value = "  -2.36%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.6271"
$r.ClearFormats()
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("E44").Value = "  -0.01%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "14.06"
$r.ClearFormats()
$ws.Range("E45").Value = "  -0.67%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "3.811"
$r.ClearFormats()
$ws.Range("E46").Value = "  +0.63%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.5852"
$r.ClearFormats()
$ws.Range("E47").Value = "  -3.06%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "2.086"
$r.ClearFormats()
$ws.Range("E48").Value = "  -0.56%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "129.54"
$r.ClearFormats()
$ws.Range("E49").Value = "  +3.26%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.228"
$r.ClearFormats()
$ws.Range("E50").Value = "  -5.55%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.07329"
$r.ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "

Write-Output "Updated cryptos sheet with new D/E values"
